# Update epexspot_prices.xlsx with the latest day of data
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column V (05-jul) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("U1").Copy()
$wsPrix.Range("V1").PasteSpecial(-4122)
$wsPrix.Range("V1").Value = "05-jul"

$vValues = @(
    41.53,
    36.29,
    37.33,
    30.42,
    30.27,
    28.73,
    28.08,
    39.34,
    17.86,
    3.34,
    0,
    -0.01,
    -0.02,
    -0.1,
    -0.11,
    -0.02,
    0.37,
    6.2,
    41.25,
    75.04000000000001,
    67.67,
    72.73999999999999,
    100.32,
    94.02
)

for ($i = 0; $i -lt $vValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 22).Value = $vValues[$i]
}

# --- Sheet "Gaz": append row 19 (2025-07-03) ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A19").NumberFormat = "@"
$wsGaz.Range("A19").Value = "2025-07-03"
$wsGaz.Range("A19").ClearFormats()
$wsGaz.Range("B19").Value = 32.85

# --- Sheet "CO2": append row 19 (2025-07-03) ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A19").NumberFormat = "@"
$wsCo2.Range("A19").Value = "2025-07-03"
$wsCo2.Range("A19").ClearFormats()
$wsCo2.Range("B19").Value = 71.81
